$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the style/format of the last existing data row (233) to use as template
$ws.Range("A233:C233").Copy()

$ws.Range("A234:C234").PasteSpecial(-4122)
$ws.Range("A234").Value = "cs"
$ws.Range("B234").Value = "lab.coil.ohm.label"
$ws.Range("C234").Value = "Odpor spirálky"

$ws.Range("A235:C235").PasteSpecial(-4122)
$ws.Range("A235").Value = "cs"
$ws.Range("B235").Value = "lab.coil.ohm.label.tooltip"
$ws.Range("C235").Value = "Odporem spirálky je v tomto případě míněn odpor v základně atomizéru."

$ws.Range("A236:C236").PasteSpecial(-4122)
$ws.Range("A236").Value = "cs"
$ws.Range("B236").Value = "lab.coil.wraps.label"
$ws.Range("C236").Value = "Počet otoček"

$ws.Range("A237:C237").PasteSpecial(-4122)
$ws.Range("A237").Value = "cs"
$ws.Range("B237").Value = "lab.coil.wraps.label.tooltip"
$ws.Range("C237").Value = "Uveďte prosím finální počet otoček (tzn. tolik, v kolika bude spirálka umístěna do atomizéru). "

$ws.Range("A238:C238").PasteSpecial(-4122)
$ws.Range("A238").Value = "cs"
$ws.Range("B238").Value = "lab.coil.code.label"
$ws.Range("C238").Value = "Kód spirálky"

$ws.Range("A239:C239").PasteSpecial(-4122)
$ws.Range("A239").Value = "cs"
$ws.Range("B239").Value = "lab.coil.code.label.tooltip"
$ws.Range("C239").Value = "Kód bude užitečný v importech, exportech a vůbec v identifikaci spirálky člověkěm. Musí být unikátní."

$ws.Range("A240:C240").PasteSpecial(-4122)
$ws.Range("A240").Value = "cs"
$ws.Range("B240").Value = "lab.wire.name.label"
$ws.Range("C240").Value = "Název drátu"

$ws.Range("A241:C241").PasteSpecial(-4122)
$ws.Range("A241").Value = "cs"
$ws.Range("B241").Value = "lab.wire.tooltip.create"
$ws.Range("C241").Value = "Vytvořit drát"

$ws.Range("A242:C242").PasteSpecial(-4122)
$ws.Range("A242").Value = "cs"
$ws.Range("B242").Value = "lab.wire.tooltip.create"
$ws.Range("C242").Value = "Vytvořit drát"

$ws.Range("A243:C243").PasteSpecial(-4122)
$ws.Range("A243").Value = "cs"
$ws.Range("B243").Value = "lab.wire.create.title"
$ws.Range("C243").Value = "Nový drát"

$ws.Range("A244:C244").PasteSpecial(-4122)
$ws.Range("A244").Value = "cs"
$ws.Range("B244").Value = "lab.wire.create.subtitle"
$ws.Range("C244").Value = "Dráty se používají k tvorbě spirálek."

$ws.Range("A245:C245").PasteSpecial(-4122)
$ws.Range("A245").Value = "cs"
$ws.Range("B245").Value = "lab.wire.description.label"
$ws.Range("C245").Value = "Popis"

$ws.Range("A246:C246").PasteSpecial(-4122)
$ws.Range("A246").Value = "cs"
$ws.Range("B246").Value = "lab.wire.name.label"
$ws.Range("C246").Value = "Jméno"

$ws.Range("A247:C247").PasteSpecial(-4122)
$ws.Range("A247").Value = "cs"
$ws.Range("B247").Value = "lab.wire.ga.label"
$ws.Range("C247").Value = "Tloušťka (GA)"

$ws.Range("A248:C248").PasteSpecial(-4122)
$ws.Range("A248").Value = "cs"
$ws.Range("B248").Value = "lab.wire.vendorId.label"
$ws.Range("C248").Value = "Výrobce"

$ws.Range("A249:C249").PasteSpecial(-4122)
$ws.Range("A249").Value = "cs"
$ws.Range("B249").Value = "lab.wire.create.submit"
$ws.Range("C249").Value = "Vytvořit drát"

$ws.Range("A250:C250").PasteSpecial(-4122)
$ws.Range("A250").Value = "cs"
$ws.Range("B250").Value = "lab.coil.create.submit"
$ws.Range("C250").Value = "Vytvořit spirálku"

$ws.Range("A251:C251").PasteSpecial(-4122)
$ws.Range("A251").Value = "cs"
$ws.Range("B251").Value = "lab.coil.create.success"
$ws.Range("C251").Value = "Spirálka byla vytvoředna."

$ws.Range("A252:C252").PasteSpecial(-4122)
$ws.Range("A252").Value = "cs"
$ws.Range("B252").Value = "lab.wire.create.success"
$ws.Range("C252").Value = "Drát byl uložen."

$excel.CutCopyMode = 0

# Update the view: scroll so row 230 is at top, and select B242 like the target
$win = $excel.ActiveWindow
$win.ScrollRow = 230
$win.ScrollColumn = 1
$ws.Range("B242").Select()
